# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.069.61"
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = "'1.835.64"
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'324.24"
$ws.Range("E5").Value = '  -3.22%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = "'0.4638"
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").Value = "'0.07862"
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("D10").Value = "'0.9618"
$ws.Range("D11").Value = "'21.94"
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").Value = "'1.830.36"
$ws.Range("E12").Value = '  +0.75%  '
$ws.Range("D13").Value = "'5.692"
$ws.Range("E13").Value = '  -2.74%  '
$ws.Range("D14").Value = "'6.921"
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("D15").Value = "'0.06851"
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").Value = "'87.48"
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = "'0.000009941"
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").Value = "'16.71"
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = "'28.082.90"
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("D22").Value = "'5.328"
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").Value = "'11.02"
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("D25").Value = "'2.034.50"
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").Value = "'154.27"
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = "'19.19"
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("D28").Value = "'5.667"
$ws.Range("E28").Value = '  -7.15%  '
$ws.Range("D29").Value = "'1.965"
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("D30").Value = "'118.14"
$ws.Range("E30").Value = '  +0.47%  '
$ws.Range("D31").Value = "'0.9368"
$ws.Range("E31").Value = '  -4.26%  '
$ws.Range("D32").Value = "'0.09230"
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("D33").Value = "'5.275"
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("D34").Value = "'1.323"
$ws.Range("E34").Value = '  -2.34%  '
$ws.Range("E35").Value = '  -4.97%  '
$ws.Range("E36").Value = '  -5.19%  '
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("D38").Value = "'1.141"
$ws.Range("E38").Value = '  -2.12%  '
$ws.Range("D39").Value = "'7.781"
$ws.Range("E39").Value = '  +2.40%  '
$ws.Range("D40").Value = "'0.5595"
$ws.Range("E40").Value = '  -2.37%  '
$ws.Range("D41").Value = "'9.892"
$ws.Range("E41").Value = '  -3.18%  '
$ws.Range("D42").Value = "'0.1764"
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("D43").Value = "'0.07262"
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = "'11.61"
$ws.Range("E44").Value = '  -1.21%  '
$ws.Range("D45").Value = "'0.5262"
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("E46").Value = '  -10.45%  '
$ws.Range("D47").Value = "'1.123"
$ws.Range("E47").Value = '  -10.00%  '
$ws.Range("D48").Value = "'1.830"
$ws.Range("E48").Value = '  -4.40%  '
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = "'1.026"
